$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46077

# Row 3
$ws.Range("C3").Value = 46077

# Row 4
$ws.Range("C4").Value = 46077

# Row 5
$ws.Range("C5").Value = 46077

# Row 6
$ws.Range("C6").Value = 46077

# Row 7
$ws.Range("C7").Value = 46077

# Row 8
$ws.Range("C8").Value = 46077

# Row 9
$ws.Range("C9").Value = 46077

# Row 10
$ws.Range("A10").Value = "A 49421-2023"
$ws.Range("B10").Value = 45211.0
$ws.Range("C10").Value = 46077
$ws.Range("F10").Value = "Kommuner"
$ws.Range("G10").Value = 1.4
$ws.Range("H10").Value = 1
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("P10").Value = 1
$ws.Range("R10").Value = "Grönfink"
$ws.Range("S10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/artfynd/A 49421-2023 artfynd.xlsx`", `"A 49421-2023`")"
$ws.Range("T10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/kartor/A 49421-2023 karta.png`", `"A 49421-2023`")"
$ws.Range("V10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomål/A 49421-2023 FSC-klagomål.docx`", `"A 49421-2023`")"
$ws.Range("W10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomålsmail/A 49421-2023 FSC-klagomål mail.docx`", `"A 49421-2023`")"
$ws.Range("X10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsyn/A 49421-2023 tillsynsbegäran.docx`", `"A 49421-2023`")"
$ws.Range("Y10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsynsmail/A 49421-2023 tillsynsbegäran mail.docx`", `"A 49421-2023`")"
$ws.Range("Z10").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/fåglar/A 49421-2023 prioriterade fågelarter.docx`", `"A 49421-2023`")"

# Row 11
$ws.Range("A11").Value = "A 24384-2023"
$ws.Range("C11").Value = 46077
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("O11").Value = 1
$ws.Range("R11").Value = "Kolflarnlav"
$ws.Range("S11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/artfynd/A 24384-2023 artfynd.xlsx`", `"A 24384-2023`")"
$ws.Range("T11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/kartor/A 24384-2023 karta.png`", `"A 24384-2023`")"
$ws.Range("V11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomål/A 24384-2023 FSC-klagomål.docx`", `"A 24384-2023`")"
$ws.Range("W11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomålsmail/A 24384-2023 FSC-klagomål mail.docx`", `"A 24384-2023`")"
$ws.Range("X11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsyn/A 24384-2023 tillsynsbegäran.docx`", `"A 24384-2023`")"
$ws.Range("Y11").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsynsmail/A 24384-2023 tillsynsbegäran mail.docx`", `"A 24384-2023`")"

# Row 12
$ws.Range("A12").Value = "A 24262-2023"
$ws.Range("C12").Value = 46077
$ws.Range("G12").Value = 2.5
$ws.Range("H12").Value = 0
$ws.Range("R12").Value = "Tibast"
$ws.Range("S12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/artfynd/A 24262-2023 artfynd.xlsx`", `"A 24262-2023`")"
$ws.Range("T12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/kartor/A 24262-2023 karta.png`", `"A 24262-2023`")"
$ws.Range("V12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomål/A 24262-2023 FSC-klagomål.docx`", `"A 24262-2023`")"
$ws.Range("W12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomålsmail/A 24262-2023 FSC-klagomål mail.docx`", `"A 24262-2023`")"
$ws.Range("X12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsyn/A 24262-2023 tillsynsbegäran.docx`", `"A 24262-2023`")"
$ws.Range("Y12").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsynsmail/A 24262-2023 tillsynsbegäran mail.docx`", `"A 24262-2023`")"

# Row 13
$ws.Range("A13").Value = "A 24368-2023"
$ws.Range("B13").Value = 45076.0
$ws.Range("C13").Value = 46077
$ws.Range("G13").Value = 3
$ws.Range("I13").Value = 1
$ws.Range("R13").Value = "Plattlummer"
$ws.Range("S13").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/artfynd/A 24368-2023 artfynd.xlsx`", `"A 24368-2023`")"
$ws.Range("T13").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/kartor/A 24368-2023 karta.png`", `"A 24368-2023`")"
$ws.Range("V13").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomål/A 24368-2023 FSC-klagomål.docx`", `"A 24368-2023`")"
$ws.Range("W13").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomålsmail/A 24368-2023 FSC-klagomål mail.docx`", `"A 24368-2023`")"
$ws.Range("X13").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsyn/A 24368-2023 tillsynsbegäran.docx`", `"A 24368-2023`")"
$ws.Range("Y13").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsynsmail/A 24368-2023 tillsynsbegäran mail.docx`", `"A 24368-2023`")"

# Row 14
$ws.Range("A14").Value = "A 24233-2023"
$ws.Range("B14").Value = 45076.0
$ws.Range("C14").Value = 46077
$ws.Range("F14").ClearContents() | Out-Null
$ws.Range("G14").Value = 3.7
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 1
$ws.Range("L14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("R14").Value = "Svartvit taggsvamp"
$ws.Range("S14").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/artfynd/A 24233-2023 artfynd.xlsx`", `"A 24233-2023`")"
$ws.Range("T14").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/kartor/A 24233-2023 karta.png`", `"A 24233-2023`")"
$ws.Range("V14").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomål/A 24233-2023 FSC-klagomål.docx`", `"A 24233-2023`")"
$ws.Range("W14").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomålsmail/A 24233-2023 FSC-klagomål mail.docx`", `"A 24233-2023`")"
$ws.Range("X14").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsyn/A 24233-2023 tillsynsbegäran.docx`", `"A 24233-2023`")"
$ws.Range("Y14").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsynsmail/A 24233-2023 tillsynsbegäran mail.docx`", `"A 24233-2023`")"
$ws.Range("Z14").ClearContents() | Out-Null

# Row 15
$ws.Range("A15").Value = "A 491-2026"
$ws.Range("B15").Value = 46029.42581018519
$ws.Range("C15").Value = 46077
$ws.Range("G15").Value = 3.2
$ws.Range("H15").Value = 1
$ws.Range("J15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("R15").Value = "Blåsippa"
$ws.Range("S15").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/artfynd/A 491-2026 artfynd.xlsx`", `"A 491-2026`")"
$ws.Range("T15").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/kartor/A 491-2026 karta.png`", `"A 491-2026`")"
$ws.Range("V15").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomål/A 491-2026 FSC-klagomål.docx`", `"A 491-2026`")"
$ws.Range("W15").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/klagomålsmail/A 491-2026 FSC-klagomål mail.docx`", `"A 491-2026`")"
$ws.Range("X15").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsyn/A 491-2026 tillsynsbegäran.docx`", `"A 491-2026`")"
$ws.Range("Y15").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0126/tillsynsmail/A 491-2026 tillsynsbegäran mail.docx`", `"A 491-2026`")"

# Row 16
$ws.Range("C16").Value = 46077

# Row 17
$ws.Range("A17").Value = "A 47173-2022"
$ws.Range("B17").Value = 44852.0
$ws.Range("C17").Value = 46077
$ws.Range("G17").Value = 1.9

# Row 18
$ws.Range("A18").Value = "A 24363-2023"
$ws.Range("C18").Value = 46077
$ws.Range("G18").Value = 2.4

# Row 19
$ws.Range("A19").Value = "A 24254-2023"
$ws.Range("C19").Value = 46077
$ws.Range("G19").Value = 1.2

# Row 20
$ws.Range("A20").Value = "A 24257-2023"
$ws.Range("B20").Value = 45076.0
$ws.Range("C20").Value = 46077
$ws.Range("G20").Value = 1

# Row 21
$ws.Range("A21").Value = "A 31246-2022"
$ws.Range("B21").Value = 44771.0
$ws.Range("C21").Value = 46077
$ws.Range("F21").ClearContents() | Out-Null
$ws.Range("G21").Value = 4.3

# Row 22
$ws.Range("A22").Value = "A 35404-2022"
$ws.Range("B22").Value = 44798.0
$ws.Range("C22").Value = 46077
$ws.Range("F22").Value = "Kommuner"
$ws.Range("G22").Value = 1.4

# Row 23
$ws.Range("A23").Value = "A 27589-2025"
$ws.Range("B23").Value = 45813.51070601852
$ws.Range("C23").Value = 46077
$ws.Range("G23").Value = 2.7

# Row 24
$ws.Range("A24").Value = "A 27561-2025"
$ws.Range("B24").Value = 45813.48378472222
$ws.Range("C24").Value = 46077
$ws.Range("G24").Value = 0.7
